# Update the "yearly" overview sheet: roll the 5 reporting periods forward
# by one year (drop 1396/12, shift 1397..1400 left, add 1401/12) and refresh
# the underlying data with the newly shifted-in figures ("read_price" refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Period headers (row 8 and row 24) ----
$periods = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$cols = @("E","F","G","H","I")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
    $ws.Range($cols[$i] + "24").Value = $periods[$i]
}

# ---- Data rows: shift existing 4 values left and append the new figure ----
function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

Set-RowValues 10 @(337, 677, 3356, 883, 1030)
Set-RowValues 12 @(108917, 195662, 236594, 561765, 660491)
Set-RowValues 13 @(4869, 13333, 9067, 16176, 16067)
Set-RowValues 15 @(597, 233, 133, 143, 119)
Set-RowValues 16 @(1578, 693, 1804, 2698, 4344)
Set-RowValues 17 @(72645, 72847, 108972, 163679, 222003)
Set-RowValues 19 @(28465, 33959, 44003, 93728, 67804)
Set-RowValues 20 @(217408, 317404, 403929, 839072, 971858)
Set-RowValues 26 @(344, 350, 364, 373, 315)
Set-RowValues 27 @(417, 440, 508, 500, 500)
